# Applies the "daily refresh" update to the HOFORS overview sheet:
#  - column C ("Förändrad") bumps from 2023-09-09 (45178) to 2023-09-10 (45179)
#    for every data row (rows 2..176)
#  - row 5 additionally gained one more signal species ("Vedticka"), so its
#    species counters (I5, Q5) increment and the species list (R5) gets the
#    new name inserted in its proper (alphabetical) place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 176

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # column C
    if ($cell.Value2() -eq 45178) {
        $cell.Value = 45179
    }
}

# Row 5 specific updates
$ws.Cells.Item(5, 9).Value = 5    # I5: Signalarter 4 -> 5
$ws.Cells.Item(5, 17).Value = 12  # Q5: Alla arter 11 -> 12

$r5 = "Knärot`r`nKoralltaggsvamp`r`nOrange taggsvamp`r`nTalltita`r`nUllticka`r`nBrandticka`r`nBronshjon`r`nSvavelriska`r`nThomsons trägnagare`r`nVedticka`r`nKopparödla`r`nBlåsippa"
$ws.Cells.Item(5, 18).Value = $r5   # R5: Artnamn list
